# Auto-generated PowerShell COM-interop script
# Applies updated market-price / profit figures to the Mandragora_Profits workbook
# (scheduled data-refresh run against the FFXIV Universalis price feed).

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 86: Filling in the Blanks | Enchanted Aurum Regis Ink
$ws.Range("H86").Value = 77772.75
$ws.Range("I86").Value = 123616.4
$ws.Range("J86").Value = 1366.6666
$ws.Range("K86").Value = 123616.4
$ws.Range("L86").Value = 1366.6666
$ws.Range("M86").Value = -122493.4
$ws.Range("N86").Value = -3612.6666
# Row 89: Ink into Antiquity (L) | Enchanted Aurum Regis Ink
$ws.Range("H89").Value = 77772.75
$ws.Range("I89").Value = 123616.4
$ws.Range("J89").Value = 1366.6666
$ws.Range("K89").Value = 618082
$ws.Range("L89").Value = 6833.333000000001
$ws.Range("M89").Value = -612466
$ws.Range("N89").Value = -18065.333
# Row 93: Spellbound | Koppranickel Index
$ws.Range("H93").Value = 46000
$ws.Range("J93").Value = 46000
$ws.Range("L93").Value = 46000
$ws.Range("N93").Value = -50992
# Row 135: For Tired Minds | Grade 1 Gemsap of Intelligence
$ws.Range("H135").Value = 744.93024
$ws.Range("I135").Value = 279.4
$ws.Range("J135").Value = 1149.7391
$ws.Range("K135").Value = 2514.6
$ws.Range("L135").Value = 10347.6519
$ws.Range("M135").Value = 20.40000000000009
$ws.Range("N135").Value = -15417.6519
# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Range("H138").Value = 1273.4407
$ws.Range("I138").Value = 1060.4062
$ws.Range("J138").Value = 1525.9259
$ws.Range("K138").Value = 3181.2186
$ws.Range("L138").Value = 4577.7777
$ws.Range("M138").Value = 1958.7814
$ws.Range("N138").Value = -14857.7777

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 1930.71
$ws.Range("I32").Value = 1704.2262
$ws.Range("J32").Value = 3119.75
$ws.Range("K32").Value = 1704.2262
$ws.Range("L32").Value = 3119.75
$ws.Range("M32").Value = -1417.2262
$ws.Range("N32").Value = -3693.75
# Row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Range("H61").Value = 2028
$ws.Range("I61").Value = 2326.7693
$ws.Range("J61").Value = 1785.25
$ws.Range("K61").Value = 2326.7693
$ws.Range("L61").Value = 1785.25
$ws.Range("M61").Value = -2114.7693
$ws.Range("N61").Value = -2209.25
# Row 74: As the Bolt Flies | Titanium Nugget
$ws.Range("H74").Value = 1226.5834
$ws.Range("I74").Value = 1056.2858
$ws.Range("J74").Value = 1465
$ws.Range("K74").Value = 1056.2858
$ws.Range("L74").Value = 1465
$ws.Range("M74").Value = -182.2858000000001
$ws.Range("N74").Value = -3213
# Row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws.Range("H77").Value = 1226.5834
$ws.Range("I77").Value = 1056.2858
$ws.Range("J77").Value = 1465
$ws.Range("K77").Value = 5281.429
$ws.Range("L77").Value = 7325
$ws.Range("M77").Value = -913.4290000000001
$ws.Range("N77").Value = -16061
# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 3882.1746
$ws.Range("I132").Value = 2553.697
$ws.Range("J132").Value = 5343.5
$ws.Range("K132").Value = 7661.091
$ws.Range("L132").Value = 16030.5
$ws.Range("M132").Value = -5131.091
$ws.Range("N132").Value = -21090.5
# Row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Range("H136").Value = 2028
$ws.Range("I136").Value = 2326.7693
$ws.Range("J136").Value = 1785.25
$ws.Range("K136").Value = 6980.3079
$ws.Range("L136").Value = 5355.75
$ws.Range("M136").Value = -4430.3079
$ws.Range("N136").Value = -10455.75

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 4247.345
$ws.Range("I134").Value = 1601.6428
$ws.Range("J134").Value = 6991.037
$ws.Range("K134").Value = 4804.928400000001
$ws.Range("L134").Value = 20973.111
$ws.Range("M134").Value = -2269.928400000001
$ws.Range("N134").Value = -26043.111

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 7: Gridania's Got Talent | Maple Lumber
$ws.Range("H7").Value = 5609.1665
$ws.Range("I7").Value = 14325.714
$ws.Range("K7").Value = 14325.714
$ws.Range("M7").Value = -14212.714
# Row 17: Say It with Spears | Feathered Harpoon
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").Value = $null
# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 5466466.5
$ws.Range("I31").Value = 1543.902
$ws.Range("J31").Value = 33337572
$ws.Range("K31").Value = 1543.902
$ws.Range("L31").Value = 33337572
$ws.Range("M31").Value = -1248.902
$ws.Range("N31").Value = -33338162
# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 5466466.5
$ws.Range("I34").Value = 1543.902
$ws.Range("J34").Value = 33337572
$ws.Range("K34").Value = 1543.902
$ws.Range("L34").Value = 33337572
$ws.Range("M34").Value = -1341.902
$ws.Range("N34").Value = -33337976

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap | Maple Syrup
$ws.Range("H5").Value = 560.26086
$ws.Range("J5").Value = 1429
$ws.Range("L5").Value = 4287
$ws.Range("N5").Value = -4511
# Row 122: Salt of the North | Northern Sea Salt
$ws.Range("H122").Value = 3188.7856
$ws.Range("I122").Value = 1099
$ws.Range("J122").Value = 3349.5386
$ws.Range("K122").Value = 9891
$ws.Range("L122").Value = 30145.8474
$ws.Range("M122").Value = -7441
$ws.Range("N122").Value = -35045.8474
# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Range("H131").Value = 1001.15
$ws.Range("I131").Value = 415.7143
$ws.Range("J131").Value = 1125.3334
$ws.Range("K131").Value = 1247.1429
$ws.Range("L131").Value = 3376.0002
$ws.Range("M131").Value = 3792.8571
$ws.Range("N131").Value = -13456.0002
# Row 135: Not-so-secret Ingredient | Royal Maple Syrup
$ws.Range("H135").Value = 560.26086
$ws.Range("J135").Value = 1429
$ws.Range("L135").Value = 12861
$ws.Range("N135").Value = -17931

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 132: On Board for Lar | Lar Ingot
$ws.Range("H132").Value = 1814092.4
$ws.Range("I132").Value = 4631542
$ws.Range("J132").Value = 2875
$ws.Range("K132").Value = 13894626
$ws.Range("L132").Value = 8625
$ws.Range("M132").Value = -13892096
$ws.Range("N132").Value = -13685

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 9: A Taste for Dalmaticae | Amateur's Dalmatica
$ws.Range("H9").Value = 10000000
$ws.Range("I9").Value = 10000000
$ws.Range("K9").Value = 10000000
$ws.Range("M9").Value = -9999860
# Row 81: Where the Dragonflies, the Net Catches | Crawler Silk
$ws.Range("H81").Value = 2253.5557
$ws.Range("I81").Value = 1762.5
$ws.Range("J81").Value = 2646.4
$ws.Range("K81").Value = 3525
$ws.Range("L81").Value = 5292.8
$ws.Range("M81").Value = -2464
$ws.Range("N81").Value = -7414.8
# Row 84: To Kill a Dragon on Nameday (L) | Crawler Silk
$ws.Range("H84").Value = 2253.5557
$ws.Range("I84").Value = 1762.5
$ws.Range("J84").Value = 2646.4
$ws.Range("K84").Value = 17625
$ws.Range("L84").Value = 26464
$ws.Range("M84").Value = -12321
$ws.Range("N84").Value = -37072
# Row 92: Modest Beginnings | Bloodhempen Culottes of Casting
$ws.Range("H92").Value = 25000
$ws.Range("J92").Value = 25000
$ws.Range("L92").Value = 25000
$ws.Range("N92").Value = -29992
# Row 122: Heavy Armoire | Dark Hempen Cloth
$ws.Range("H122").Value = 2709.5454
$ws.Range("I122").Value = 1225
$ws.Range("J122").Value = 6668.3335
$ws.Range("K122").Value = 3675
$ws.Range("L122").Value = 20005.0005
$ws.Range("M122").Value = -1225
$ws.Range("N122").Value = -24905.0005
# Row 126: A Polished Purchase | Snow Linen
$ws.Range("H126").Value = 1688.7073
$ws.Range("I126").Value = 1860.6875
$ws.Range("J126").Value = 1077.2222
$ws.Range("K126").Value = 5582.0625
$ws.Range("L126").Value = 3231.6666
$ws.Range("M126").Value = -3112.0625
$ws.Range("N126").Value = -8171.6666
# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 1856.0488
$ws.Range("I132").Value = 2244.611
$ws.Range("J132").Value = 1551.9565
$ws.Range("K132").Value = 6733.833
$ws.Range("L132").Value = 4655.8695
$ws.Range("M132").Value = -4203.833
$ws.Range("N132").Value = -9715.869500000001
# Row 136: Weaving the Envelope | Sarcenet Cloth
$ws.Range("H136").Value = 29413018
$ws.Range("I136").Value = 45455350
$ws.Range("J136").Value = 2066.6667
$ws.Range("K136").Value = 136366050
$ws.Range("L136").Value = 6200.000100000001
$ws.Range("M136").Value = -136363500
$ws.Range("N136").Value = -11300.0001
